$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text would otherwise be auto-detected as a number by Excel.
# These are written with a leading apostrophe (forces text entry) and then
# restored to the "Normal" style so no stray number-format style sticks around.
$forceTextCells = @{
    "D4" = '1.002'
    "D5" = '246.76'
    "D6" = '1.002'
    "D8" = '0.2975'
    "D9" = '0.06828'
    "D10" = '19.50'
    "D11" = '106.35'
    "D12" = '0.07758'
    "D14" = '5.427'
    "D15" = '0.7160'
    "D16" = '283.40'
    "D18" = '0.000007762'
    "D19" = '13.23'
    "D20" = '1.002'
    "D22" = '5.500'
    "D23" = '1.003'
    "D24" = '6.605'
    "D25" = '9.892'
    "D26" = '169.79'
    "D27" = '20.44'
    "D28" = '2.223'
    "D29" = '0.1064'
    "D30" = '1.443'
    "D31" = '4.715'
    "D32" = '1.594'
    "D33" = '4.464'
    "D34" = '0.05025'
    "D35" = '0.7620'
    "D36" = '1.163'
    "D37" = '2.737'
    "D38" = '0.02048'
    "D39" = '2.708'
    "D40" = '2.180'
    "D41" = '6.416'
    "D42" = '0.4543'
    "D43" = '109.93'
    "D44" = '0.8843'
    "D45" = '72.31'
    "D46" = '1.001'
    "D47" = '7.654'
    "D48" = '0.2649'
    "D49" = '965.09'
    "D50" = '9.435'
    "D51" = '0.1269'
}

# Cells whose target text is unambiguously non-numeric (coin names, URLs,
# padded percent strings) -- plain assignment keeps them as text already.
$plainCells = @{
    "D2" = '31.184.94'
    "E2" = '  +1.76%  '
    "D3" = '1.957.28'
    "E3" = '  +1.02%  '
    "E4" = '  +0.26%  '
    "E5" = '  +0.17%  '
    "E6" = '  +0.18%  '
    "E7" = '  +1.26%  '
    "B8" = 'Cardano'
    "C8" = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
    "E8" = '  +2.05%  '
    "B9" = 'Dogecoin'
    "C9" = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
    "E9" = '  +0.76%  '
    "B10" = 'Solana'
    "C10" = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
    "E10" = '  +0.73%  '
    "B11" = 'Litecoin'
    "C11" = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
    "E11" = '  -5.92%  '
    "B12" = 'TRON'
    "C12" = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    "E12" = '  +1.99%  '
    "B13" = 'WrappedEther'
    "C13" = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    "D13" = '1.937.59'
    "E13" = '  +0.00%  '
    "B14" = 'Polkadot'
    "C14" = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
    "E14" = '  -2.08%  '
    "B15" = 'Polygon'
    "C15" = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    "E15" = '  +5.08%  '
    "B16" = 'BitcoinCash'
    "C16" = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    "E16" = '  -4.30%  '
    "B17" = 'WrappedBTC'
    "C17" = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    "D17" = '31.061.74'
    "E17" = '  +1.29%  '
    "B18" = 'ShibaInu'
    "C18" = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
    "E18" = '  +1.36%  '
    "B19" = 'Avalanche'
    "C19" = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
    "E19" = '  +0.52%  '
    "B20" = 'Dai'
    "C20" = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
    "E20" = '  +0.19%  '
    "B21" = 'WrappedliquidstakedEther2.0'
    "C21" = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    "D21" = '2.190.35'
    "E21" = '  -0.06%  '
    "B22" = 'Uniswap'
    "C22" = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
    "E22" = '  -1.08%  '
    "B23" = 'BinanceUSD'
    "C23" = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
    "E23" = '  +0.35%  '
    "B24" = 'Chainlink'
    "C24" = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
    "E24" = '  +1.63%  '
    "B25" = 'Cosmos'
    "C25" = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
    "E25" = '  +2.82%  '
    "B26" = 'Monero'
    "C26" = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    "E26" = '  +1.07%  '
    "B27" = 'EthereumClassic'
    "C27" = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    "E27" = '  +0.59%  '
    "B28" = 'LidoDAOToken'
    "C28" = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    "E28" = '  +5.65%  '
    "B29" = 'Stellar'
    "C29" = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    "E29" = '  -0.63%  '
    "B30" = 'Toncoin'
    "C30" = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
    "E30" = '  +1.16%  '
    "B31" = 'Filecoin'
    "C31" = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    "E31" = '  +14.87%  '
    "B32" = 'PancakeSwap'
    "C32" = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    "E32" = '  +0.92%  '
    "E33" = '  +6.79%  '
    "E34" = '  +0.52%  '
    "E35" = '  +1.11%  '
    "E36" = '  +0.67%  '
    "E37" = '  +0.51%  '
    "E38" = '  -1.06%  '
    "E39" = '  +0.57%  '
    "E40" = '  +7.73%  '
    "E41" = '  +9.54%  '
    "E42" = '  +2.23%  '
    "E43" = '  -0.26%  '
    "E44" = '  +1.31%  '
    "E45" = '  +1.91%  '
    "E46" = '  -0.03%  '
    "E47" = '  +3.64%  '
    "B48" = 'WOONetwork'
    "C48" = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
    "E48" = '  +4.01%  '
    "B49" = 'Maker'
    "C49" = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    "E49" = '  +12.98%  '
    "E50" = '  +1.28%  '
    "B51" = 'Algorand'
    "C51" = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    "E51" = '  +3.45%  '
}

foreach ($addr in $forceTextCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $forceTextCells[$addr]
    $cell.Style = "Normal"
}

foreach ($addr in $plainCells.Keys) {
    $ws.Range($addr).Value = $plainCells[$addr]
}

